$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.759.11"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.436.18"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.01"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.91"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.439.26"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -4.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "4.038.54"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.61"
$ws.Range("E15").Value = "  -3.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").Value = "  -6.06%  "
$ws.Range("D17").Value = "64.798.74"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "3.431.75"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.87"
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.81"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.95"
$ws.Range("E22").Value = "  -4.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.11"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -4.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.81"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.09"
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.01"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.21"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.01"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.17"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").Value = "2.901.23"
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0746"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.66"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.26"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.54"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.88"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0316"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.777"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.89"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.09"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "315.59"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.51"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("E51").Value = "  -4.24%  "
